$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 867 (pushes the existing 867:908 block down to 868:909)
$ws.Rows.Item(867).Insert()

# Column A holds dates stored as plain text (e.g. "2026/02/26") in this sheet,
# so force Text format before writing to avoid Excel auto-converting the
# literal into a real date serial value.
$ws.Range("A867").NumberFormat = "@"
$ws.Range("A867").Value = "2026/02/27"
$ws.Range("B867").Value = "金"
$ws.Range("C867").Value = 4
$ws.Range("D867").Value = 30
